# Edit workbook to match the target diff:
# - Remove the "Indiana State Fairgrounds & Event Center" row (old row 14),
#   shifting rows 15-34 up to rows 14-33.
# - Update the standalone A4 rank value (6 -> 7).
# - Correct a handful of rank / metric values that deviate from a pure shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standalone rank correction on row 4 (Avon Community Heritage Festival)
$ws.Cells.Item(4, 1).Value = 7

# Delete the "Indiana State Fairgrounds & Event Center" row entirely; this
# shifts all subsequent rows up by one and shrinks the used range to A1:E33.
$ws.Rows.Item(14).Delete()

# After the shift, a handful of rank/metric values need correcting so the
# final values match the target state exactly.

# Row 15: Indianapolis Grapevine -> rank 24
$ws.Cells.Item(15, 1).Value = 24

# Row 18: Indianapolis Zoo -> rank 25, ratings total 14982
$ws.Cells.Item(18, 1).Value = 25
$ws.Cells.Item(18, 5).Value = 14982

# Row 24: MasterWorks Festival -> rank 14
$ws.Cells.Item(24, 1).Value = 14

# Row 28: Spirit & Place Festival -> rank 6
$ws.Cells.Item(28, 1).Value = 6

# Row 33: Waterman's Family Farm -> rank 13, rating 4.4, ratings total 669
$ws.Cells.Item(33, 1).Value = 13
$ws.Cells.Item(33, 4).Value = 4.4
$ws.Cells.Item(33, 5).Value = 669
